$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.792.06'
$ws.Range('E2').Value = '  -0.17%  '

# Row 3
$ws.Range('D3').Value = '1.767.80'
$ws.Range('E3').Value = '  -2.60%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.97%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '338.68'
$ws.Range('E5').Value = '  +0.40%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +1.15%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3787'
$ws.Range('E7').Value = '  -3.85%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3384'
$ws.Range('E8').Value = '  -3.07%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '46.12'
$ws.Range('E9').Value = '  -5.54%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.130'
$ws.Range('E10').Value = '  -6.53%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07303'
$ws.Range('E11').Value = '  -3.66%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.15'
$ws.Range('E12').Value = '  +4.07%  '

# Row 13
$ws.Range('E13').Value = '  +0.98%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.287'
$ws.Range('E14').Value = '  -3.87%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.324'
$ws.Range('E15').Value = '  +1.62%  '

# Row 16
$ws.Range('D16').Value = '1.771.00'
$ws.Range('E16').Value = '  -2.42%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001062'
$ws.Range('E17').Value = '  -4.09%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06642'
$ws.Range('E18').Value = '  -0.57%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '81.55'
$ws.Range('E19').Value = '  -4.41%  '

# Row 20
$ws.Range('E20').Value = '  +1.10%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.26'
$ws.Range('E21').Value = '  -3.42%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.367'
$ws.Range('E22').Value = '  -3.51%  '

# Row 23
$ws.Range('D23').Value = '27.813.98'
$ws.Range('E23').Value = '  -0.23%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.85'
$ws.Range('E24').Value = '  -8.03%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.393'
$ws.Range('E25').Value = '  -1.02%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.487'
$ws.Range('E26').Value = '  -0.93%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.31'
$ws.Range('E27').Value = '  -5.31%  '

# Row 28
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '151.62'
$ws.Range('E28').Value = '  -2.42%  '

# Row 29
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.351'
$ws.Range('E29').Value = '  -8.60%  '

# Row 30
$ws.Range('D30').Value = '1.972.25'
$ws.Range('E30').Value = '  -2.33%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '133.09'
$ws.Range('E31').Value = '  -2.01%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.047'
$ws.Range('E32').Value = '  +0.60%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.915'
$ws.Range('E33').Value = '  -3.36%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08805'
$ws.Range('E34').Value = '  -0.13%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.47'
$ws.Range('E35').Value = '  -6.21%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02373'
$ws.Range('E36').Value = '  -1.99%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6721'
$ws.Range('E37').Value = '  -2.93%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06304'
$ws.Range('E38').Value = '  -3.48%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.194'
$ws.Range('E39').Value = '  -6.54%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2134'
$ws.Range('E40').Value = '  -4.51%  '

# Row 41
$ws.Range('E41').Value = '  -7.88%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.216'
$ws.Range('E42').Value = '  -3.74%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.117'
$ws.Range('E43').Value = '  -5.17%  '

# Row 44
$ws.Range('E44').Value = '  +1.06%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.91'
$ws.Range('E45').Value = '  -5.69%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6141'
$ws.Range('E46').Value = '  -6.29%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.851'
$ws.Range('E47').Value = '  -0.21%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '131.95'
$ws.Range('E48').Value = '  -0.50%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.035'
$ws.Range('E49').Value = '  -6.03%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07317'
$ws.Range('E50').Value = '  +1.28%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.189'
$ws.Range('E51').Value = '  +2.40%  '
